$d = $word.ActiveDocument

# 1. Willard and Spackman's occupational therapy -> new text (no leading space)
$d.Content.Find.Execute(
    " Willard and Spackman's occupational therapy", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "A new approach to defining functional ability in ankylosing spondylitis: the development of the Bath Ankylosing Spondylitis Functional Index.",
    2) | Out-Null

# 2. Pedretti's Occupational Therapy-E-Book... -> new text (keeps leading space)
$d.Content.Find.Execute(
    " Pedretti's Occupational Therapy-E-Book: Practice Skills for Physical Dysfunction", $true, $false, $false, $false, $false,
    $true, 1, $false,
    " The European Spondylarthropathy Study Group preliminary criteria for the classification of spondylarthropathy",
    2) | Out-Null

# 3. Concepts of occupational therapy -> new text (no leading space)
$d.Content.Find.Execute(
    " Concepts of occupational therapy", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "The development of Assessment of Spondyloarthritis International Society (ASAS) classification criteria for axial spondyloarthritis (part II): validation and final selection",
    2) | Out-Null

# 4. Conceptual foundations of occupational therapy practice -> new text (no leading space)
$d.Content.Find.Execute(
    " Conceptual foundations of occupational therapy practice", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Defining spinal mobility in ankylosing spondylitis (AS). The Bath AS Metrology Index.",
    2) | Out-Null

# 5. Occupational therapy for children -> new text (no leading space)
$d.Content.Find.Execute(
    " Occupational therapy for children", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Ankylosing spondylitis",
    2) | Out-Null

# 6. A model of human occupation: Theory and application -> new text (no leading space)
$d.Content.Find.Execute(
    " A model of human occupation: Theory and application", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Criteria of the classification of spondylarthropathies",
    2) | Out-Null

# 7. Adult norms for the Box and Block Test of manual dexterity -> new text
$d.Content.Find.Execute(
    "Adult norms for the Box and Block Test of manual dexterity", $true, $false, $false, $false, $false,
    $true, 1, $false,
    "Treatment of active ankylosing spondylitis with infliximab: a randomised controlled multicentre trial",
    2) | Out-Null

# 8. Remove the final paragraph entirely (Relative contributions of neural mechanisms...)
$last = $d.Paragraphs($d.Paragraphs.Count)
$last.Range.Delete()
